$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to clean snake_case names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the "de/del/la/las/los/el" connector words in place names ---
$ws.Range("B8").Value = "Marqués De Comillas"
$ws.Range("B12").Value = "San Cristóbal De Las Casas"
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A27").Value = "Estado De México"
$ws.Range("B29").Value = "San Felipe Del Progreso"
$ws.Range("B30").Value = "Tlalnepantla De Baz"
$ws.Range("B31").Value = "Valle De Chalco Solidaridad"
$ws.Range("B33").Value = "Apaseo El Alto"
$ws.Range("B38").Value = "Acapulco De Juárez"
$ws.Range("B40").Value = "Atoyac De Álvarez"
$ws.Range("B41").Value = "Ayutla De Los Libres"
$ws.Range("B44").Value = "Iguala De La Independencia"
$ws.Range("B45").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B50").Value = "Zihuatanejo De Azueta"
$ws.Range("B53").Value = "Atotonilco De Tula"
$ws.Range("B57").Value = "Pachuca De Soto"
$ws.Range("B61").Value = "Cuautitlán De García Barragán"
$ws.Range("A65").Value = "Michoacán De Ocampo"
$ws.Range("B87").Value = "Teotitlán De Flores Magón"
$ws.Range("B92").Value = "Ixcamilpa De Guerrero"
$ws.Range("B98").Value = "Landa De Matamoros"
$ws.Range("A113").Value = "Veracruz De Ignacio De La Llave"

# --- Grand total label: all-caps -> title case ---
$ws.Range("A131").Value = "Total"

# --- Drop the trailing footnote / source rows (133-137) ---
$ws.Rows.Item(137).Delete()
$ws.Rows.Item(136).Delete()
$ws.Rows.Item(135).Delete()
$ws.Rows.Item(134).Delete()
$ws.Rows.Item(133).Delete()
